# Insert a new data row before the existing row 404 (shifts old rows
# 404..480 down to 405..481) and populate the new row 404 with the
# new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 404 (and everything below it) down by one row.
$ws.Rows.Item(404).EntireRow.Insert()

# Fill in the new row 404 with the latest observation. Columns that
# are constant for every row in this sheet (market/category metadata)
# are copied straight across; the price-related columns hold the new
# values.
$ws.Range("A404").Value = 8
$ws.Range("B404").Value = "Terminal La Palmera de La Serena"
$ws.Range("C404").Value = "Coquimbo"
$ws.Range("D404").Value = 44637
$ws.Range("D404").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E404").Value = 4
$ws.Range("F404").Value = 100112043
$ws.Range("G404").Value = "Pepino ensalada"
$ws.Range("H404").Value = "Sin especificar"
$ws.Range("I404").Value = "Primera"
$ws.Range("J404").Value = 680
$ws.Range("K404").Value = 18000
$ws.Range("L404").Value = 19000
$ws.Range("M404").Value = 18500
$ws.Range("N404").Value = "$/caja 60 unidades"
$ws.Range("O404").Value = "Región de Arica y Parinacota"
$ws.Range("P404").Value = 308
$ws.Range("Q404").Value = 60
$ws.Range("R404").Value = "Hortaliza"
